# Natmi following Dr Hou advice
# Update LR-pair (Fn1-Itgb3) results sheet with recomputed expression /
# specificity values (ligand- and receptor-expressing cell counts changed
# from 1 to 3 per group, with all downstream totals / specificities
# recomputed accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.11671699999999
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.214110666666667
$ws.Range("N2").Value = 21.642332
$ws.Range("O2").Value = 0.4688823795981188
$ws.Range("P2").Value = 0.4688823795981188
$ws.Range("Q2").Value = 195.0616577848938
$ws.Range("R2").Value = 1755.554920064044
$ws.Range("S2").Value = 0.03327277618216672
$ws.Range("T2").Value = 0.03327277618216672

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.11671699999999
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.110350666666666
$ws.Range("N3").Value = 21.331052
$ws.Range("O3").Value = 0.4621384803214003
$ws.Range("P3").Value = 0.4621384803214003
$ws.Range("Q3").Value = 192.2561009329204
$ws.Range("R3").Value = 1730.304908396284
$ws.Range("S3").Value = 0.03279421639618871
$ws.Range("T3").Value = 0.03279421639618871

# Row 4 (ECs -> sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.11671699999999
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.061296333333333
$ws.Range("N4").Value = 3.183889
$ws.Range("O4").Value = 0.06897914008048092
$ws.Range("P4").Value = 0.06897914008048092
$ws.Range("Q4").Value = 28.69629144137922
$ws.Range("R4").Value = 258.266622972413
$ws.Range("S4").Value = 0.00489488961198186
$ws.Range("T4").Value = 0.004894889611981861

# Row 5 (FAPs -> ECs)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.214110666666667
$ws.Range("N5").Value = 21.642332
$ws.Range("O5").Value = 0.4688823795981188
$ws.Range("P5").Value = 0.4688823795981188
$ws.Range("Q5").Value = 2492.955543607409
$ws.Range("R5").Value = 22436.59989246668
$ws.Range("S5").Value = 0.4252376032096085
$ws.Range("T5").Value = 0.4252376032096085

# Row 6 (FAPs -> FAPs)
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.110350666666666
$ws.Range("N6").Value = 21.331052
$ws.Range("O6").Value = 0.4621384803214003
$ws.Range("P6").Value = 0.4621384803214003
$ws.Range("Q6").Value = 2457.099555370369
$ws.Range("R6").Value = 22113.89599833332
$ws.Range("S6").Value = 0.4191214434017335
$ws.Range("T6").Value = 0.4191214434017335

# Row 7 (FAPs -> sCs)
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.061296333333333
$ws.Range("N7").Value = 3.183889
$ws.Range("O7").Value = 0.06897914008048092
$ws.Range("P7").Value = 0.06897914008048092
$ws.Range("Q7").Value = 366.7485432152436
$ws.Range("R7").Value = 3300.736888937192
$ws.Range("S7").Value = 0.06255838452369351
$ws.Range("T7").Value = 0.06255838452369351

# Row 8 (sCs -> ECs)
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.214110666666667
$ws.Range("N8").Value = 21.642332
$ws.Range("O8").Value = 0.4688823795981188
$ws.Range("P8").Value = 0.4688823795981188
$ws.Range("Q8").Value = 60.80585352174578
$ws.Range("R8").Value = 547.2526816957121
$ws.Range("S8").Value = 0.01037200020634369
$ws.Range("T8").Value = 0.0103720002063437

# Row 9 (sCs -> FAPs)
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.110350666666666
$ws.Range("N9").Value = 21.331052
$ws.Range("O9").Value = 0.4621384803214003
$ws.Range("P9").Value = 0.4621384803214003
$ws.Range("Q9").Value = 59.93128759769245
$ws.Range("R9").Value = 539.3815883792321
$ws.Range("S9").Value = 0.01022282052347816
$ws.Range("T9").Value = 0.01022282052347816

# Row 10 (sCs -> sCs)
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.061296333333333
$ws.Range("N10").Value = 3.183889
$ws.Range("O10").Value = 0.06897914008048092
$ws.Range("P10").Value = 0.06897914008048092
$ws.Range("Q10").Value = 8.945389441558222
$ws.Range("R10").Value = 80.50850497402401
$ws.Range("S10").Value = 0.001525865944805551
$ws.Range("T10").Value = 0.001525865944805552
